# Scheduled Sheets runner: refresh FFXIV Diabolos market-board derived values
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) across all crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1228.6
$ws.Range("I4").Value = 326.33334
$ws.Range("K4").Value = 326.33334
$ws.Range("M4").Value = -212.33334
$ws.Range("H9").Value = 35715220
$ws.Range("I9").Value = 41667600
$ws.Range("K9").Value = 41667600
$ws.Range("M9").Value = -41667431
$ws.Range("H51").Value = 5099.5654
$ws.Range("J51").Value = 5552.353
$ws.Range("L51").Value = 5552.353
$ws.Range("N51").Value = -6520.353
$ws.Range("H62").Value = 96509.44500000001
$ws.Range("I62").Value = 63758.8
$ws.Range("K62").Value = 63758.8
$ws.Range("M62").Value = -63134.8
$ws.Range("H65").Value = 96509.44500000001
$ws.Range("I65").Value = 63758.8
$ws.Range("K65").Value = 318794
$ws.Range("M65").Value = -315674
$ws.Range("H107").Value = 999.3333
$ws.Range("I107").Value = 999.3333
$ws.Range("K107").Value = 999.3333
$ws.Range("M107").Value = 920.6667
$ws.Range("H112").Value = 1740.1177
$ws.Range("J112").Value = 1817.625
$ws.Range("L112").Value = 5452.875
$ws.Range("N112").Value = -7668.875
$ws.Range("H113").Value = 66670610
$ws.Range("J113").Value = 5143
$ws.Range("L113").Value = 5143
$ws.Range("N113").Value = -11651
$ws.Range("H132").Value = 4945.25
$ws.Range("I132").Value = 4565.75
$ws.Range("J132").Value = 7222.25
$ws.Range("K132").Value = 13697.25
$ws.Range("L132").Value = 21666.75
$ws.Range("M132").Value = -11167.25
$ws.Range("N132").Value = -26726.75
$ws.Range("H141").Value = 2998.75
$ws.Range("I141").Value = 2998.75
$ws.Range("K141").Value = 8996.25
$ws.Range("M141").Value = -3816.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1659.6666
$ws.Range("I2").Value = 1659.6666
$ws.Range("K2").Value = 1659.6666
$ws.Range("M2").Value = -1546.6666
$ws.Range("H32").Value = 2367.963
$ws.Range("I32").Value = 2401.3845
$ws.Range("K32").Value = 2401.3845
$ws.Range("M32").Value = -2114.3845
$ws.Range("H61").Value = 55558810
$ws.Range("I61").Value = 66669570
$ws.Range("K61").Value = 66669570
$ws.Range("M61").Value = -66669358
$ws.Range("H116").Value = 1659.6666
$ws.Range("I116").Value = 1659.6666
$ws.Range("K116").Value = 1659.6666
$ws.Range("M116").Value = 634.3334
$ws.Range("H122").Value = 10419687
$ws.Range("I122").Value = 13336040
$ws.Range("K122").Value = 40008120
$ws.Range("M122").Value = -40005670
$ws.Range("H136").Value = 55558810
$ws.Range("I136").Value = 66669570
$ws.Range("K136").Value = 200008710
$ws.Range("M136").Value = -200006160

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1659.6666
$ws.Range("I3").Value = 1659.6666
$ws.Range("K3").Value = 1659.6666
$ws.Range("M3").Value = -1545.6666
$ws.Range("H13").Value = 76000
$ws.Range("I13").Value = 75000
$ws.Range("J13").Value = 77000
$ws.Range("K13").Value = 75000
$ws.Range("L13").Value = 77000
$ws.Range("M13").Value = -74832
$ws.Range("N13").Value = -77336
$ws.Range("H82").Value = 47269.082
$ws.Range("I82").Value = 15064.25
$ws.Range("K82").Value = 15064.25
$ws.Range("M82").Value = -14681.25
$ws.Range("H85").Value = 47269.082
$ws.Range("I85").Value = 15064.25
$ws.Range("K85").Value = 15064.25
$ws.Range("M85").Value = -13738.25
$ws.Range("H134").Value = 1805.7333
$ws.Range("I134").Value = 1513.6
$ws.Range("K134").Value = 4540.799999999999
$ws.Range("M134").Value = -2005.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2787.8594
$ws.Range("I31").Value = 1496.7391
$ws.Range("J31").Value = 3512.1462
$ws.Range("K31").Value = 1496.7391
$ws.Range("L31").Value = 3512.1462
$ws.Range("M31").Value = -1201.7391
$ws.Range("N31").Value = -4102.1462
$ws.Range("H34").Value = 2787.8594
$ws.Range("I34").Value = 1496.7391
$ws.Range("J34").Value = 3512.1462
$ws.Range("K34").Value = 1496.7391
$ws.Range("L34").Value = 3512.1462
$ws.Range("M34").Value = -1294.7391
$ws.Range("N34").Value = -3916.1462
$ws.Range("H107").Value = 1350.3684
$ws.Range("I107").Value = 1324.4482
$ws.Range("K107").Value = 1324.4482
$ws.Range("M107").Value = 595.5518

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 412.51724
$ws.Range("I2").Value = 831.38464
$ws.Range("J2").Value = 72.1875
$ws.Range("K2").Value = 4988.307839999999
$ws.Range("L2").Value = 433.125
$ws.Range("M2").Value = -4875.307839999999
$ws.Range("N2").Value = -659.125
$ws.Range("H38").Value = 139.05556
$ws.Range("I38").Value = 160.2
$ws.Range("K38").Value = 480.6
$ws.Range("M38").Value = -133.6
$ws.Range("H63").Value = 17006
$ws.Range("I63").Value = 17006
$ws.Range("K63").Value = 51018
$ws.Range("M63").Value = -50269
$ws.Range("H66").Value = 17006
$ws.Range("I66").Value = 17006
$ws.Range("K66").Value = 153054
$ws.Range("M66").Value = -149310
$ws.Range("H112").Value = 2461.25
$ws.Range("I112").Value = 2461.25
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 7383.75
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -6275.75
$ws.Range("N112").ClearContents()
$ws.Range("H122").Value = 836.2
$ws.Range("J122").Value = 892.8570999999999
$ws.Range("L122").Value = 8035.7139
$ws.Range("N122").Value = -12935.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1605.25
$ws.Range("I107").Value = 559.6667
$ws.Range("J107").Value = 2650.8333
$ws.Range("K107").Value = 559.6667
$ws.Range("L107").Value = 2650.8333
$ws.Range("M107").Value = 1360.3333
$ws.Range("N107").Value = -6490.8333
$ws.Range("H113").Value = 1984.5834
$ws.Range("I113").Value = 820.3333
$ws.Range("K113").Value = 820.3333
$ws.Range("M113").Value = 1349.6667
$ws.Range("H122").Value = 17859234
$ws.Range("I122").Value = 2039.0416
$ws.Range("J122").Value = 125002400
$ws.Range("K122").Value = 6117.1248
$ws.Range("L122").Value = 375007200
$ws.Range("M122").Value = -3667.1248
$ws.Range("N122").Value = -375012100
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 75000
$ws.Range("J130").Value = 75000
$ws.Range("L130").Value = 75000
$ws.Range("N130").Value = -85040
$ws.Range("H132").Value = 3853.4412
$ws.Range("I132").Value = 3636.0454
$ws.Range("K132").Value = 10908.1362
$ws.Range("M132").Value = -8378.136200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H117").Value = 35000
$ws.Range("J117").Value = 35000
$ws.Range("L117").Value = 35000
$ws.Range("N117").Value = -44178
$ws.Range("H122").Value = 3650.3809
$ws.Range("I122").Value = 2589.077
$ws.Range("J122").Value = 5375
$ws.Range("K122").Value = 7767.231000000001
$ws.Range("L122").Value = 16125
$ws.Range("M122").Value = -5317.231000000001
$ws.Range("N122").Value = -21025
$ws.Range("H136").Value = 2193.9473
$ws.Range("I136").Value = 1992.2858
$ws.Range("J136").Value = 2758.6
$ws.Range("K136").Value = 5976.857400000001
$ws.Range("L136").Value = 8275.799999999999
$ws.Range("M136").Value = -3426.857400000001
$ws.Range("N136").Value = -13375.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4958.5
$ws.Range("I81").Value = 1253.0625
$ws.Range("K81").Value = 2506.125
$ws.Range("M81").Value = -1445.125
$ws.Range("H84").Value = 4958.5
$ws.Range("I84").Value = 1253.0625
$ws.Range("K84").Value = 12530.625
$ws.Range("M84").Value = -7226.625
$ws.Range("H126").Value = 1462.0952
$ws.Range("I126").Value = 1415.2
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 4245.6
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -1775.6
$ws.Range("N126").Value = -12140
$ws.Range("H132").Value = 2560.3901
$ws.Range("I132").Value = 1833.8438
$ws.Range("K132").Value = 5501.5314
$ws.Range("M132").Value = -2971.5314

Write-Output "Applied scheduled Sheets update."